$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 133
$ws.Range("A133").Value = 7
$ws.Range("B133").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C133").Value = "Ñuble"
$ws.Range("D133").Value = 44911
$ws.Range("D133").NumberFormat = $ws.Range("D132").NumberFormat
$ws.Range("E133").Value = 16
$ws.Range("F133").Value = "Fruta"
$ws.Range("G133").Value = 100103
$ws.Range("H133").Value = "Frutos de hueso (carozo)"
$ws.Range("I133").Value = 100103001
$ws.Range("J133").Value = "Cereza"
$ws.Range("K133").Value = "Lapins"
$ws.Range("L133").Value = "Primera"
$ws.Range("M133").Value = 160
$ws.Range("N133").Value = 4500
$ws.Range("O133").Value = 5000
$ws.Range("P133").Value = 4750
$ws.Range("Q133").Value = "$/bandeja 10 kilos"
$ws.Range("R133").Value = "Provincia de Curicó"
$ws.Range("S133").Value = 475
$ws.Range("T133").Value = 10

# New row 134
$ws.Range("A134").Value = 7
$ws.Range("B134").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C134").Value = "Ñuble"
$ws.Range("D134").Value = 44911
$ws.Range("D134").NumberFormat = $ws.Range("D132").NumberFormat
$ws.Range("E134").Value = 16
$ws.Range("F134").Value = "Fruta"
$ws.Range("G134").Value = 100103
$ws.Range("H134").Value = "Frutos de hueso (carozo)"
$ws.Range("I134").Value = 100103001
$ws.Range("J134").Value = "Cereza"
$ws.Range("K134").Value = "Lapins"
$ws.Range("L134").Value = "Segunda"
$ws.Range("M134").Value = 80
$ws.Range("N134").Value = 3000
$ws.Range("O134").Value = 3000
$ws.Range("P134").Value = 3000
$ws.Range("Q134").Value = "$/bandeja 10 kilos"
$ws.Range("R134").Value = "Provincia de Curicó"
$ws.Range("S134").Value = 300
$ws.Range("T134").Value = 10
